$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.934.83"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.309.66"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "103.10"
$ws.Range("E5").Value = "  +5.52%  "
$ws.Range("D6").Value = "271.44"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.609"
$ws.Range("E9").Value = "  -2.95%  "
$ws.Range("D10").Value = "45.71"
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("D11").Value = "0.0937"
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("D12").Value = "7.95"
$ws.Range("E12").Value = "  -2.30%  "
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("D14").Value = "15.93"
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("D15").Value = "2.648.76"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").Value = "0.862"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("D17").Value = "2.293.63"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "43.826.46"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").Value = "6.27"
$ws.Range("E20").Value = "  -2.49%  "
$ws.Range("D21").Value = "72.44"
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("D22").Value = "2.47"
$ws.Range("E22").Value = "  +8.00%  "
$ws.Range("D23").Value = "233.85"
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("D24").Value = "2.93"
$ws.Range("E24").Value = "  +15.64%  "
$ws.Range("D25").Value = "9.20"
$ws.Range("E25").Value = "  -2.81%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "11.27"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "2.29"
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "38.62"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "177.43"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("D32").Value = "21.91"
$ws.Range("E32").Value = "  -2.50%  "
$ws.Range("D33").Value = "0.0896"
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").Value = "5.49"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").Value = "4.80"
$ws.Range("E36").Value = "  +7.77%  "
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").Value = "0.0354"
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("D39").Value = "3.60"
$ws.Range("E39").Value = "  +6.07%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "2.36"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.236"
$ws.Range("E41").Value = "  -1.78%  "
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("D43").Value = "12.26"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").Value = "64.83"
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("D45").Value = "8.86"
$ws.Range("E45").Value = "  -3.89%  "
$ws.Range("D46").Value = "5.26"
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").Value = "98.89"
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").Value = "0.447"
$ws.Range("E50").Value = "  +6.65%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "1.52"
$ws.Range("E51").Value = "  +10.61%  "
